# Apply Kraken_Profits market-price/profit updates scraped by the scheduled runner.
# Each block updates the H/I/J/K/L/M/N (price & profit) columns for one leve row
# on the given sheet; some rows gain or lose a trailing profit cell entirely.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 18
$ws.Range("H18").Value = 23200
$ws.Range("I18").Value = 23200
$ws.Range("K18").Value = 23200
$ws.Range("M18").Value = -22916

# ALC row 20
$ws.Range("H20").Value = 15724.75
$ws.Range("I20").Value = 7633
$ws.Range("K20").Value = 7633
$ws.Range("M20").Value = -7403

# ALC row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null

# ALC row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null

# ALC row 26
$ws.Range("H26").Value = 50000
$ws.Range("J26").Value = 50000
$ws.Range("L26").Value = 50000
$ws.Range("N26").Value = -50688

# ALC row 31
$ws.Range("H31").Value = 1366.6666
$ws.Range("I31").Value = 1366.6666
$ws.Range("K31").Value = 4099.9998
$ws.Range("M31").Value = -3869.9998

# ALC row 35
$ws.Range("H35").Value = 15724.75
$ws.Range("I35").Value = 7633
$ws.Range("K35").Value = 7633
$ws.Range("M35").Value = -7254

# ALC row 43
$ws.Range("H43").Value = 2886
$ws.Range("J43").Value = 2886
$ws.Range("L43").Value = 2886
$ws.Range("N43").Value = -3024

# ALC row 55
$ws.Range("H55").Value = 264.84616
$ws.Range("I55").Value = 123.666664
$ws.Range("J55").Value = 582.5
$ws.Range("K55").Value = 123.666664
$ws.Range("L55").Value = 582.5
$ws.Range("M55").Value = 90.333336
$ws.Range("N55").Value = -1010.5

# ALC row 88
$ws.Range("H88").Value = 3174.8333
$ws.Range("I88").Value = 3512.25
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 3512.25
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -3106.25
$ws.Range("N88").Value = -3312

# ALC row 91
$ws.Range("H91").Value = 3174.8333
$ws.Range("I91").Value = 3512.25
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 3512.25
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -2108.25
$ws.Range("N91").Value = -5308

# ALC row 97
$ws.Range("H97").Value = 3893.75
$ws.Range("J97").Value = 3893.75
$ws.Range("L97").Value = 11681.25
$ws.Range("N97").Value = -12673.25

# ALC row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null

$ws = $wb.Worksheets.Item("ARM")
# ARM row 30
$ws.Range("H30").Value = 597.4
$ws.Range("I30").Value = 597.4
$ws.Range("K30").Value = 597.4
$ws.Range("M30").Value = -447.4

# ARM row 88
$ws.Range("H88").Value = 535.5
$ws.Range("I88").Value = 464.33334
$ws.Range("J88").Value = 578.2
$ws.Range("K88").Value = 464.33334
$ws.Range("L88").Value = 578.2
$ws.Range("M88").Value = -58.33334000000002
$ws.Range("N88").Value = -1390.2

# ARM row 91
$ws.Range("H91").Value = 535.5
$ws.Range("I91").Value = 464.33334
$ws.Range("J91").Value = 578.2
$ws.Range("K91").Value = 464.33334
$ws.Range("L91").Value = 578.2
$ws.Range("M91").Value = 939.66666
$ws.Range("N91").Value = -3386.2

# ARM row 110
$ws.Range("H110").Value = 1702.5264
$ws.Range("I110").Value = 1099.875
$ws.Range("K110").Value = 1099.875
$ws.Range("M110").Value = 945.125

# ARM row 122
$ws.Range("H122").Value = 8197.066000000001
$ws.Range("J122").Value = 4598.6
$ws.Range("L122").Value = 13795.8
$ws.Range("N122").Value = -18695.8

# ARM row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").Value = $null

$ws = $wb.Worksheets.Item("BSM")
# BSM row 22
$ws.Range("H22").Value = 334448.66
$ws.Range("I22").Value = 501173
$ws.Range("K22").Value = 501173
$ws.Range("M22").Value = -501000

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58
$ws.Range("H58").Value = 3487.75
$ws.Range("I58").Value = 3795.8
$ws.Range("K58").Value = 3795.8
$ws.Range("M58").Value = -3592.8

# CRP row 60
$ws.Range("H60").Value = 40000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 40000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 40000
$ws.Range("M60").Value = $null
$ws.Range("N60").Value = -41022

# CRP row 136
$ws.Range("H136").Value = 3487.75
$ws.Range("I136").Value = 3795.8
$ws.Range("K136").Value = 11387.4
$ws.Range("M136").Value = -8837.400000000001

$ws = $wb.Worksheets.Item("CUL")
# CUL row 37
$ws.Range("H37").Value = 59999
$ws.Range("J37").Value = 59999
$ws.Range("L37").Value = 179997
$ws.Range("N37").Value = -180221

# CUL row 54
$ws.Range("H54").Value = 2000
$ws.Range("J54").Value = 2000
$ws.Range("L54").Value = 6000
$ws.Range("N54").Value = -7118

# CUL row 55
$ws.Range("H55").Value = 1990
$ws.Range("I55").Value = 1366.6666
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 4099.9998
$ws.Range("L55").Value = 7500
$ws.Range("M55").Value = -3922.9998
$ws.Range("N55").Value = -7854

# CUL row 139
$ws.Range("H139").Value = 1000
$ws.Range("I139").Value = 1000
$ws.Range("K139").Value = 3000
$ws.Range("M139").Value = 2140

$ws = $wb.Worksheets.Item("GSM")
# GSM row 5
$ws.Range("H5").Value = 2835.3333
$ws.Range("I5").Value = 1253
$ws.Range("K5").Value = 1253
$ws.Range("M5").Value = -1141

# GSM row 102
$ws.Range("H102").Value = 2722.5386
$ws.Range("I102").Value = 2363.4546
$ws.Range("K102").Value = 2363.4546
$ws.Range("M102").Value = -741.4546

# GSM row 122
$ws.Range("H122").Value = 7929.5713
$ws.Range("J122").Value = 8000
$ws.Range("L122").Value = 24000
$ws.Range("N122").Value = -28900

# GSM row 132
$ws.Range("H132").Value = 2671.4736
$ws.Range("I132").Value = 2009.5714
$ws.Range("K132").Value = 6028.7142
$ws.Range("M132").Value = -3498.7142

$ws = $wb.Worksheets.Item("LTW")
# LTW row 100
$ws.Range("H100").Value = 6284.5713
$ws.Range("I100").Value = 3064.7778
$ws.Range("K100").Value = 3064.7778
$ws.Range("M100").Value = -2523.7778

# LTW row 122
$ws.Range("H122").Value = 6500
$ws.Range("I122").Value = 6500
$ws.Range("K122").Value = 19500
$ws.Range("M122").Value = -17050

# LTW row 136
$ws.Range("H136").Value = 1189.8
$ws.Range("I136").Value = 1112.25
$ws.Range("K136").Value = 3336.75
$ws.Range("M136").Value = -786.75

$ws = $wb.Worksheets.Item("WVR")
# WVR row 24
$ws.Range("H24").Value = 2500499.5
$ws.Range("I24").Value = 2500499.5
$ws.Range("K24").Value = 2500499.5
$ws.Range("M24").Value = -2500269.5

# WVR row 81
$ws.Range("H81").Value = 200
$ws.Range("I81").Value = 200
$ws.Range("K81").Value = 400
$ws.Range("M81").Value = 661

# WVR row 84
$ws.Range("H84").Value = 200
$ws.Range("I84").Value = 200
$ws.Range("K84").Value = 2000
$ws.Range("M84").Value = 3304
